$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "customer_id" column before the existing "group_id" column.
$ws.Columns("B:B").Insert()
$ws.Range("B1").Value = "customer_id"
$ws.Range("B2:B5").Value = 1

# The insert shifts the old hyperlink cell from E5 to F5; the engine does not
# renumber the hyperlink's own range automatically, so re-home it explicitly
# (re-adding keeps the same cell value/shared-string and relationship id).
$ws.Range("E5").Hyperlinks.Delete() | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), "https://t.me/testInteractTool") | Out-Null
$ws.Range("F5").Style = "Hyperlink"

# Resize the new customer_id column and the (now shifted) group_id column.
$ws.Columns("B:B").ColumnWidth = 28.15
$ws.Columns("C:C").ColumnWidth = 31.33

# Append duplicate group rows (2,3,4) for a second customer_id = 2.
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = -1001429987581
$ws.Range("D6").Value = "test"
$ws.Range("E6").Value = "private"
$ws.Range("F6").Value = "VTvg_eT6s7Rz-AIj"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = -1001170310837
$ws.Range("D7").Value = "TestKDbot"
$ws.Range("E7").Value = "private"
$ws.Range("F7").Value = "RcGGtdG60NynCrJK"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 2
$ws.Range("C8").Value = -1001159430667
$ws.Range("D8").Value = "Test Tool"
$ws.Range("E8").Value = "public"
$ws.Range("F8").Value = "https://t.me/testInteractTool"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://t.me/testInteractTool") | Out-Null
$ws.Range("F8").Style = "Hyperlink"

# Match the workbook's final selection state.
$ws.Range("C12").Select() | Out-Null
